$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 1.41

# Row 3
$ws.Range("G3").Value = 1.24

# Row 6
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 2.68
$ws.Range("H6").Value = 1.59
$ws.Range("I6").Value = 5.2
$ws.Range("J6").Value = 3.25
$ws.Range("P6").Value = 2

# Row 9
$ws.Range("Q9").Value = 2.46

# Row 11
$ws.Range("N11").Value = 3.45
$ws.Range("O11").Value = 1.38
$ws.Range("P11").Value = 1.83
$ws.Range("Q11").Value = 2.14
$ws.Range("S11").Value = 3.95
$ws.Range("T11").Value = 1.86
$ws.Range("X11").Value = 12.5
$ws.Range("Y11").Value = 10.5
$ws.Range("AB11").Value = 10.5
$ws.Range("AH11").Value = 18.5
$ws.Range("AI11").Value = 50
$ws.Range("AK11").Value = 36
$ws.Range("AL11").Value = 50
$ws.Range("AN11").Value = 34
$ws.Range("AO11").Value = 34
